$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read TS data from excel to DataTable / Process DataTable into dict /
# Create a TimeSeries and add data to it.
# -> B2 now derives its value from C2 instead of holding a literal.
$ws.Range("B2").Formula = "=C2"

# Restore the active selection to B3 (as left by the author after editing).
$ws.Range("B3").Select()

# Touch the page setup so the sheet carries explicit print settings.
$ws.PageSetup.Orientation = 1
